$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$commitSha = "38fcabc93432b923d908fd6169f3569ee1047fef"
$newFile = "6ac94990-807c-4f21-9691-96402d440f76.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFile"

# ---------------------------------------------------------------------------
# Sheet "Overview" - add row 3
# ---------------------------------------------------------------------------
$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = "e2e\$newFile"
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-25 04:38:49"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), $newFileUrl, $null, $null, "e2e\$newFile")

$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - add row 3
# ---------------------------------------------------------------------------
$ws2.Range("A3").Value = $newFile
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "False"
$ws2.Range("G3").Value = "6ac94990-807c-4f21-9691-96402d440f76.025e5b97701041d5fa1fcd01f052f0ab7919d58f.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-25 04:38:44"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

$ws2.Hyperlinks.Add($ws2.Range("A3"), $newFileUrl, $null, $null, $newFile)

$tbl2 = $ws2.ListObjects.Item(1)
$tbl2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" - add row 3
# ---------------------------------------------------------------------------
$ws3.Range("A3").Value = $newFile
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "False"
$ws3.Range("G3").Value = "6ac94990-807c-4f21-9691-96402d440f76.025e5b97701041d5fa1fcd01f052f0ab7919d58f.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-25 04:38:49"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

$ws3.Hyperlinks.Add($ws3.Range("A3"), $newFileUrl, $null, $null, $newFile)

$tbl3 = $ws3.ListObjects.Item(1)
$tbl3.Resize($ws3.Range("A1:P3"))

Write-Host "Overview dims: $($ws1.UsedRange.Address())"
Write-Host "zh-cn dims: $($ws2.UsedRange.Address())"
Write-Host "de-de dims: $($ws3.UsedRange.Address())"
